$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section: CARBON and NITROGEN (row 32 header, rows 33-36 data) ---
$ws.Cells.Item(32, 1).Value = "CARBON and NITROGEN"
$ws.Cells.Item(32, 1).Font.Underline = 2

$ws.Cells.Item(33, 1).Value = "Logistic Regression"
$ws.Cells.Item(33, 2).Value = 82.4

$ws.Cells.Item(34, 1).Value = "Decision Tree"
$ws.Cells.Item(34, 2).Value = 91.4

$ws.Cells.Item(35, 1).Value = "kNN"
$ws.Cells.Item(35, 2).Value = 90.7

$ws.Cells.Item(36, 1).Value = "Random Forest Classifier"
$ws.Cells.Item(36, 2).Value = 95.8

# --- Section: C, Si, N (row 38 header, rows 39-42 data) ---
$ws.Cells.Item(38, 1).Value = "C, Si, N"
$ws.Cells.Item(38, 1).Font.Underline = 2

$ws.Cells.Item(39, 1).Value = "kNN"
$ws.Cells.Item(39, 2).Value = 87.1

$ws.Cells.Item(40, 1).Value = "Logistic Regression"
$ws.Cells.Item(40, 2).Value = 84.4

$ws.Cells.Item(41, 1).Value = "Decision Tree"
$ws.Cells.Item(41, 2).Value = 95.7

$ws.Cells.Item(42, 1).Value = "Random Forest Classifier"
$ws.Cells.Item(42, 2).Value = 96.9

# --- View state: scroll/selection to match final saved view ---
[void]$ws.Range("C40").Select()
$excel.ActiveWindow.ScrollRow = 26

Write-Host "done"
